$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values.
#    Typed in the same order the original author entered them so that the
#    shared-string table ends up built in the same sequence as the target
#    workbook (row2 A/B, row3 A/B/C, row2 C, row4 A/B/C, row5 A/B/C).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "takeDamage() function"
$ws.Range("B2").Value = "Instead of decreasing health by damage, health is set to be current damage value."
$ws.Range("A3").Value = "resetDamage() function"
$ws.Range("B3").Value = "Instead of setting damage to 0, it is setting to damage itself."
$ws.Range("C3").Value = "Instead od damage = damage, I've put damage = 0."
$ws.Range("C2").Value = "Instead of health = damage, I've put health -= damage."
$ws.Range("A4").Value = "getBoost() function"
$ws.Range("B4").Value = "Instead of returning boost, it returns number."
$ws.Range("C4").Value = "Instead of return number, I've out return boost."
$ws.Range("A5").Value = "findNumberInHand() function"
$ws.Range("B5").Value = "True and false in return statements are swapped."
$ws.Range("C5").Value = "When we found card, function now returns true, not false."

# ---------------------------------------------------------------------------
# 2. Header row touch-up: B1 used to carry its own near-duplicate style;
#    make it match A1/C1 (single-property change on a single cell keeps the
#    existing shared style instead of growing the style table).
# ---------------------------------------------------------------------------
$ws.Range("B1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Body formatting.
#    Build each target look once on an out-of-the-way helper cell, then
#    copy/paste-special the resulting format onto the real ranges. Doing it
#    this way (one fully-resolved style per paste) avoids the runtime
#    generating extra throw-away style entries that a sequence of
#    property-by-property edits on a multi-cell range would otherwise leave
#    behind. Each contiguous block is pasted separately because this
#    runtime's PasteSpecial only honours the first area of a multi-area
#    destination range.
# ---------------------------------------------------------------------------
$helper = $ws.Range("Z100")

# Rows 2, 3 and 5: left aligned, vertically centered, wrapped text.
$helper.HorizontalAlignment = -4131
$helper.VerticalAlignment = -4108
$helper.WrapText = $true
$helper.Copy()
$ws.Range("A2:C3").PasteSpecial(-4122)
$ws.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$helper.Clear()

# Row 4: same alignment, but with the font explicitly re-applied (as in the
# authored workbook, this row ends up using a distinct-but-visually-equal
# font entry).
$helper.Font.Name = "Arial"
$helper.Font.Size = 10
$helper.HorizontalAlignment = -4131
$helper.VerticalAlignment = -4108
$helper.WrapText = $true
$helper.Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$helper.Clear()

# ---------------------------------------------------------------------------
# 4. Row heights: rows 2, 3 and 5 wrap onto two lines.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 25.5
$ws.Rows(3).RowHeight = 25.5
$ws.Rows(5).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 5. Column widths (closest achievable values; the runtime quantizes column
#    widths to a coarser grid than native Excel).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 35.42
$ws.Columns("B").ColumnWidth = 44.59
$ws.Columns("C").ColumnWidth = 40.59

Write-Host "edit complete"
